# Commit: "add employee Test case."
# Adds a new "CreateEmployee" worksheet (after "GroupDetails") that holds a
# small employee-record test fixture (Location/FullName/employeeID/... plus
# 3 sample rows), with a couple of alignment styles and hyperlinked emails.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new worksheet at the end of the tab strip and name it.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "CreateEmployee"

# ---------------------------------------------------------------------
# 2. Header row (row 1) — reuses a couple of existing shared strings
#    ("Location") and introduces the rest of the field-name strings.
# ---------------------------------------------------------------------
$ws3.Cells.Item(1, 1).Value = "Location"
$ws3.Cells.Item(1, 2).Value = "FullName"
$ws3.Cells.Item(1, 3).Value = "employeeID"
$ws3.Cells.Item(1, 4).Value = "supervisor"
$ws3.Cells.Item(1, 5).Value = "businessUnit"
$ws3.Cells.Item(1, 6).Value = "department"
$ws3.Cells.Item(1, 7).Value = "address"
$ws3.Cells.Item(1, 8).Value = "city"
$ws3.Cells.Item(1, 9).Value = "state"
$ws3.Cells.Item(1, 10).Value = "postalCode"
$ws3.Cells.Item(1, 11).Value = "country"
$ws3.Cells.Item(1, 12).Value = "phone"
$ws3.Cells.Item(1, 13).Value = "email"

# ---------------------------------------------------------------------
# 3. First data record (row 2). Field order here matches the order the
#    source fixture's values were authored in.
# ---------------------------------------------------------------------
$ws3.Cells.Item(2, 6).Value = "IT"
$ws3.Cells.Item(2, 7).Value = "45, GB Thosmson Road"
$ws3.Cells.Item(2, 2).Value = "Alex Starch"
$ws3.Cells.Item(2, 8).Value = "NewyorkCity"
$ws3.Cells.Item(2, 9).Value = "NY"
$ws3.Cells.Item(2, 11).Value = "usa"
$ws3.Cells.Item(2, 13).Value = "alexstarch@testmail.com"
$ws3.Cells.Item(2, 12).Value = "866-782-4964"
$ws3.Cells.Item(2, 4).Value = "John Smith"
$ws3.Cells.Item(2, 1).Value = "Dallas"
$ws3.Cells.Item(2, 3).Value = 4567
$ws3.Cells.Item(2, 10).Value = 10019

# ---------------------------------------------------------------------
# 4. Rows 3 and 4 repeat the same record (row 3 omits FullName).
# ---------------------------------------------------------------------
$ws3.Cells.Item(3, 1).Value = "Dallas"
$ws3.Cells.Item(3, 3).Value = 4567
$ws3.Cells.Item(3, 4).Value = "John Smith"
$ws3.Cells.Item(3, 6).Value = "IT"
$ws3.Cells.Item(3, 7).Value = "45, GB Thosmson Road"
$ws3.Cells.Item(3, 8).Value = "NewyorkCity"
$ws3.Cells.Item(3, 9).Value = "NY"
$ws3.Cells.Item(3, 10).Value = 10019
$ws3.Cells.Item(3, 11).Value = "usa"
$ws3.Cells.Item(3, 12).Value = "866-782-4964"
$ws3.Cells.Item(3, 13).Value = "alexstarch@testmail.com"

$ws3.Cells.Item(4, 1).Value = "Dallas"
$ws3.Cells.Item(4, 2).Value = "Alex Starch"
$ws3.Cells.Item(4, 3).Value = 4567
$ws3.Cells.Item(4, 4).Value = "John Smith"
$ws3.Cells.Item(4, 6).Value = "IT"
$ws3.Cells.Item(4, 7).Value = "45, GB Thosmson Road"
$ws3.Cells.Item(4, 8).Value = "NewyorkCity"
$ws3.Cells.Item(4, 9).Value = "NY"
$ws3.Cells.Item(4, 10).Value = 10019
$ws3.Cells.Item(4, 11).Value = "usa"
$ws3.Cells.Item(4, 12).Value = "866-782-4964"
$ws3.Cells.Item(4, 13).Value = "alexstarch@testmail.com"

# ---------------------------------------------------------------------
# 5. Cell styles. Apply LEFT alignment before CENTER alignment so the
#    generated stylesheet allocates the xf records in that order
#    (matches the authored workbook: xf#2 = left, xf#3 = center).
# ---------------------------------------------------------------------
$xlLeft = -4131
$xlCenter = -4108

$ws3.Range("C2").HorizontalAlignment = $xlLeft
$ws3.Range("C3").HorizontalAlignment = $xlLeft
$ws3.Range("C4").HorizontalAlignment = $xlLeft

$ws3.Range("G1").HorizontalAlignment = $xlCenter
$ws3.Range("G2").HorizontalAlignment = $xlCenter
$ws3.Range("G3").HorizontalAlignment = $xlCenter
$ws3.Range("G4").HorizontalAlignment = $xlCenter
$ws3.Range("J1").HorizontalAlignment = $xlCenter
$ws3.Range("J2").HorizontalAlignment = $xlCenter
$ws3.Range("J3").HorizontalAlignment = $xlCenter
$ws3.Range("J4").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------
# 6. Hyperlink the email column (adds the Hyperlink cell style too).
# ---------------------------------------------------------------------
$ws3.Hyperlinks.Add($ws3.Range("M2"), "mailto:alexstarch@testmail.com")
$ws3.Hyperlinks.Add($ws3.Range("M3"), "mailto:alexstarch@testmail.com")
$ws3.Hyperlinks.Add($ws3.Range("M4"), "mailto:alexstarch@testmail.com")

# ---------------------------------------------------------------------
# 7. Column widths for the new sheet.
# ---------------------------------------------------------------------
$ws3.Columns.Item(1).ColumnWidth = 12.17
$ws3.Columns.Item(2).ColumnWidth = 14.75
$ws3.Columns.Item(3).ColumnWidth = 16.17
$ws3.Columns.Item(4).ColumnWidth = 13.33
$ws3.Columns.Item(5).ColumnWidth = 14.5
$ws3.Columns.Item(6).ColumnWidth = 11.17
$ws3.Columns.Item(7).ColumnWidth = 26.58
$ws3.Columns.Item(8).ColumnWidth = 15.75
$ws3.Columns.Item(9).ColumnWidth = 10.33
$ws3.Columns.Item(10).ColumnWidth = 16.5
$ws3.Columns.Item(11).ColumnWidth = 11.33
$ws3.Columns.Item(12).ColumnWidth = 14.33
$ws3.Columns.Item(13).ColumnWidth = 12.5

# ---------------------------------------------------------------------
# 8. Selection / active cell, matching the authored file. GroupDetails'
#    own selection moves too (it is no longer the active tab).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("GroupDetails").Range("C1").Select()
$ws3.Activate()
$ws3.Range("H15").Select()
